# "Add files via upload" - populate the Repository (F) and URL (H) columns
# for rows 2-7 with links to the author's GitHub project, and turn the
# first Repository cell (F2) into a live hyperlink.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$repositoryUrl = "https://github.com/majedunl/Cornelia-Sorabji.github.io"
$sourceUrl     = "https://github.com/majedunl/Cornelia-Sorabji-"

# Column H = URL, Column F = Repository.
# Write H first so the shared-string table picks up the URL text (index 73)
# before the repository text (index 74), matching the order the values were
# entered in.
for ($row = 2; $row -le 7; $row++) {
    $ws.Cells.Item($row, 8).Value = $sourceUrl
    $ws.Cells.Item($row, 6).Value = $repositoryUrl
}

# Make F2 a real hyperlink back to the source URL, keeping the existing
# built-in "Hyperlink" cell style.
$ws.Hyperlinks.Add($ws.Range("F2"), $sourceUrl)
$ws.Range("F2").Style = "Hyperlink"

# Leave the selection where the editor last left it.
$ws.Range("F7").Select()
